# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" / "Office" colour scheme (used by the notes master)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" colour scheme (used by the slide master,
#                            i.e. it is the theme that actually paints every slide)
# The target edit swaps the two themes' content: theme1.xml ends up holding the
# "Integral"/Red Violet colours and theme2.xml ends up holding the plain "Office Theme"
# colours (everything else - fonts, format scheme - is identical between the two parts).
#
# The only theme this host exposes through the PowerPoint object model is the one that is
# actually applied to the slides/slide master (ppt/theme/theme2.xml), reached through
# Slide.ThemeColorScheme (12 slots, in the canonical DrawingML order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink). We drive the visible, effective result of the swap by
# repainting those 12 slots with the colours the "Office Theme" originally used.

function Color-RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order (1-based) exposed by Slide.ThemeColorScheme:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$tcs.Item(1).RGB  = Color-RGB 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = Color-RGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = Color-RGB 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = Color-RGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = Color-RGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = Color-RGB 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = Color-RGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = Color-RGB 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = Color-RGB 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = Color-RGB 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = Color-RGB 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = Color-RGB 0x95 0x4F 0x72   # folHlink 954F72
